# Daily attendance processing - 2025-10-12 07:18:00
# For every "Recorded By" (column G) cell on the Session Analysis Results
# sheet whose comma-separated list of recorders starts with the automated
# "System" entry, move that leading entry to the back of the list so the
# real (human/service) recorder is listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*"
        if ($parts.Count -gt 1 -and $parts[0].Trim().ToLower() -eq "system") {
            $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
            $cell.Value = [string]::Join(", ", $rotated)
        }
    }
}
